$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update existing timestamps (rows 2-7) to the new scrape time
$newTs = "2025-09-21 18:22:52"
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTs
}

# 2. Insert new row 5 (GAS) - push WordPress/電子秤/Jotform down
$ws.Rows.Item(5).Insert()
$ws.Cells.Item(5, 1).Value = $newTs
$ws.Cells.Item(5, 2).Value = "【GAS開発者募集】Amazon広告管理SaaSのMVP開発"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5397812"
$ws.Cells.Item(5, 7).Value = 88
$ws.Cells.Item(5, 8).Value = "◆開発 ◇管理"

# 3. Insert new row 8 (Gemini) - before Jotform (now at row 8), push it to row 9
$ws.Rows.Item(8).Insert()
$ws.Cells.Item(8, 1).Value = $newTs
$ws.Cells.Item(8, 2).Value = "Geminiで旅行のしおりのHTMLを生成するプロンプトの作成"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5397817"
$ws.Cells.Item(8, 7).Value = 10

# 4. Rebuild hyperlinks for column F (rows 2-9), since Insert() does not shift them.
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le 9; $r++) {
    $url = $ws.Cells.Item($r, 6).Value2
    $ws.Hyperlinks.Add($ws.Range("F" + $r), $url, "", "", $url)
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
}
